# Apply the cell-value edits described by the diff to Sheet1 ("Joel L.")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "[-, 'MEC-3B-Tec. Fundição']"
$ws.Range("F2").Value = "-"

# Row 4
$ws.Range("B4").Value = "-"

# Row 6
$ws.Range("C6").Value = "['MCT-1A-Tecnologia dos Materiais.', -]"
$ws.Range("E6").Value = "['MEC-3B-Tec. Fundição', -, -, -]"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "['MEC-3B-Tec. Fundição', -, -, -]"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "[-, -, -, 'MEC-3B-Tec. Fundição']"
$ws.Range("F8").Value = "-"

# Row 10
$ws.Range("B10").Value = "-"
$ws.Range("D10").Value = "[-, -, -, 'MEC-3A-Tec. Fundição']"

# Row 11
$ws.Range("D11").Value = "[-, -, -, 'MEC-3A-Tec. Fundição']"

# Row 12
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "[-, -, -, 'MEC-3A-Tec. Fundição']"

# Row 14
$ws.Range("B14").Value = "-"

# Row 16
$ws.Range("B16").Value = "-"
$ws.Range("D16").Value = "['MEC-3A-Tec. Fundição', -, -, -]"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "['MEC-2NB-Fundição', -, -, -]"
$ws.Range("F18").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "[Rogério-Retífica-2NB, -, -, -]"
$ws.Range("D19").Value = "-"
$ws.Range("F19").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"

# Row 20
$ws.Range("C20").Value = "['MEC-2NB-Fundição', -, -, -]"
$ws.Range("D20").Value = "-"
$ws.Range("F20").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "['MEC-2NB-Fundição', -, -, -]"
$ws.Range("D21").Value = "-"
$ws.Range("F21").Value = "[-, -, -, 'MEC-1NA-T. M. Metalicos']"
